$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.961.40'
$ws.Range("E2").Value = '  -2.79%  '

$ws.Range("D3").Value = '3.475.08'
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("E5").Value = '  -1.70%  '

$ws.Range("D6").Value = '173.74'
$ws.Range("E6").Value = '  -3.66%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  -2.96%  '

$ws.Range("D9").Value = '3.475.62'
$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("E10").Value = '  -6.29%  '

$ws.Range("D11").Value = '6.86'
$ws.Range("E11").Value = '  -1.65%  '

$ws.Range("D13").Value = '4.078.98'
$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("E14").Value = '  +0.16%  '

$ws.Range("D15").Value = '30.00'
$ws.Range("E15").Value = '  -6.53%  '

$ws.Range("D16").Value = '66.070.89'
$ws.Range("E16").Value = '  -2.67%  '

$ws.Range("E17").Value = '  -3.32%  '

$ws.Range("D18").Value = '3.475.92'
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("E19").Value = '  -4.04%  '

$ws.Range("D20").Value = '13.91'
$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").Value = '366.39'
$ws.Range("E21").Value = '  -6.56%  '

$ws.Range("D22").Value = '7.75'
$ws.Range("E22").Value = '  -1.93%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.27%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '72.41'
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("E27").Value = '  -7.80%  '

$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  +1.51%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("D30").Value = '23.94'
$ws.Range("E30").Value = '  +1.97%  '

$ws.Range("E31").Value = '  -3.12%  '

$ws.Range("E32").Value = '  -5.71%  '

$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("E34").Value = '  -2.85%  '

$ws.Range("E35").Value = '  -7.48%  '

$ws.Range("E36").Value = '  -1.68%  '

$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").Value = '29.37'
$ws.Range("E37").Value = '  +12.62%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '159.15'
$ws.Range("E38").Value = '  -1.83%  '

$ws.Range("D39").Value = '0.889'
$ws.Range("E39").Value = '  +0.30%  '

$ws.Range("D40").Value = '2.815.46'
$ws.Range("E40").Value = '  +3.26%  '

$ws.Range("E41").Value = '  -5.72%  '

$ws.Range("D42").Value = '2.57'
$ws.Range("E42").Value = '  -7.13%  '

$ws.Range("E43").Value = '  -3.88%  '

$ws.Range("D44").Value = '6.40'
$ws.Range("E44").Value = '  -4.26%  '

$ws.Range("E45").Value = '  -4.96%  '

$ws.Range("D46").Value = '39.90'
$ws.Range("E46").Value = '  -3.24%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '24.11'
$ws.Range("E47").Value = '  -7.68%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0289'
$ws.Range("E48").Value = '  -2.89%  '

$ws.Range("D49").Value = '310.57'
$ws.Range("E49").Value = '  -5.52%  '

$ws.Range("E50").Value = '  -2.41%  '

$ws.Range("E51").Value = '  -2.12%  '

